$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Urenverantwoording")

# New rows of hour-log data to append (rows 79-94)
# Columns: A=Taak, B=Persoon, C=Uren, D=Datum, E=Beschrijving
$rows = @(
    @("Software Architecture Solution", "Arco",   6,   42668, "STD's maken"),
    @("Software Architecture Solution", "Marten", 6,   42668, "STD's maken"),
    @("Software Architecture Solution", "Brent",  6,   42668, "Objectlijst maken"),
    @("Software Architecture Solution", "Rene",   6,   42668, "Objectlijst maken"),
    @("Opmaak & Indeling",              "Arco",   1,   42668, "Opzet onderzoek"),
    @("Opmaak & Indeling",              "Marten", 1,   42668, "Opzet onderzoek"),
    @("Opmaak & Indeling",              "Brent",  1,   42668, "Opzet onderzoek"),
    @("Opmaak & Indeling",              "Rene",   1,   42668, "Opzet onderzoek"),
    @("Software Architecture Solution", "Arco",   4.5, 42669, "STD's afmaken"),
    @("Software Architecture Solution", "Marten", 2.5, 42669, "STD's afmaken"),
    @("Software Architecture Solution", "Brent",  4.5, 42669, "STD's afmaken"),
    @("Software Architecture Solution", "Rene",   4.5, 42669, "STD's afmaken"),
    @("Inhoud",                         "Brent",  2,   42669, "Bedenken welke RTOS'en er gebruikt gaan worden"),
    @("Inhoud",                         "Rene",   2,   42669, "Bedenken welke RTOS'en er gebruikt gaan worden"),
    @("Administratie",                  "Arco",   1,   42670, "Administratie van de afgelopen tijd samengevoegd tot totale tijd"),
    @("Software Architecture Solution", "Arco",   0.5, 42670, "STD's samenvoegen")
)

$startRow = 79
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]
    $ws.Cells.Item($r, 1).Value = $data[0]
    $ws.Cells.Item($r, 2).Value = $data[1]
    $ws.Cells.Item($r, 3).Value = $data[2]
    $ws.Cells.Item($r, 4).Value = [double]$data[3]
    $ws.Cells.Item($r, 5).Value = $data[4]
}

# Match the date-cell style used by the rest of column D (same as D78)
$ws.Range("D79:D94").NumberFormat = "[$-F800]dddd\,\ mmmm\ dd\,\ yyyy"

$excel.CalculateFullRebuild()
